$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it sat between " không" and
#    " đồng nhất thông qua màn hình Serial."). It will be re-created at
#    the end of the newly inserted paragraph below.
# ----------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ----------------------------------------------------------------------
# 2. Split the paragraph right after "...màn hình Serial." (i.e. right
#    before its paragraph mark) so that sentence stays in its own
#    paragraph and a brand-new paragraph follows it.
# ----------------------------------------------------------------------
$targetPara = $d.Paragraphs.Item(5)
$splitPoint = $targetPara.Range.End - 1
$d.Range($splitPoint, $splitPoint).InsertParagraphAfter()

# ----------------------------------------------------------------------
# 3. Fill the freshly created (empty) paragraph with the new sentence,
#    keeping "15s" bold.
# ----------------------------------------------------------------------
$newPara = $d.Paragraphs.Item(6)
$newPara.Range.InsertAfter("Từ lúc tạo request cho tới lúc nhận response của get https là khoảng 15s.")

$newParaStart = $newPara.Range.Start
$newParaText = $newPara.Range.Text
$boldOffset = $newParaText.IndexOf("15s")
$boldStart = $newParaStart + $boldOffset
$boldEnd = $boldStart + 3
$d.Range($boldStart, $boldEnd).Font.Bold = 1

# ----------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark at the end of this new paragraph,
#    right before its paragraph mark (mirrors the original placement
#    style, just moved to the new location).
#
#    NOTE: placing a *collapsed* bookmark range exactly at
#    `paragraph.Range.End - 1` (the position immediately before a
#    paragraph mark) is mishandled by this COM host and silently resets
#    to document position 0. Work around it by temporarily inserting a
#    throw-away character after the target spot (so the spot is no
#    longer "last position before the paragraph mark"), adding the
#    bookmark there, then deleting the throw-away character again - the
#    bookmark stays put.
# ----------------------------------------------------------------------
$bookmarkPos = $newPara.Range.End - 1
$d.Range($bookmarkPos, $bookmarkPos).InsertAfter("~")
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
$d.Range($bookmarkPos, $bookmarkPos + 1).Delete()

# ----------------------------------------------------------------------
# 5. Add a new, otherwise empty, paragraph (containing a single space)
#    right after the sentence paragraph.
# ----------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.InsertAfter(" ")
